$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Efnb2'
$ws.Cells.Item(2, 3).Value = 'Grm1'
$ws.Cells.Item(2, 4).Value = 'FAPs'
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 53.0169945
$ws.Cells.Item(2, 8).Value = 106.033989
$ws.Cells.Item(2, 9).Value = 0.7442421144210264
$ws.Cells.Item(2, 10).Value = 0.7041603619966643
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.003113666666666667
$ws.Cells.Item(2, 14).Value = 0.009341
$ws.Cells.Item(2, 15).Value = 0.005310654226435927
$ws.Cells.Item(2, 16).Value = 0.005310654226435927
$ws.Cells.Item(2, 17).Value = 0.1650772485415
$ws.Cells.Item(2, 18).Value = 0.9904634912490001
$ws.Cells.Item(2, 19).Value = 0.003952412530441635
$ws.Cells.Item(2, 20).Value = 0.003739552202526238

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Efnb2'
$ws.Cells.Item(3, 3).Value = 'Grm1'
$ws.Cells.Item(3, 4).Value = 'Neutrophils'
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 53.0169945
$ws.Cells.Item(3, 8).Value = 106.033989
$ws.Cells.Item(3, 9).Value = 0.7442421144210264
$ws.Cells.Item(3, 10).Value = 0.7041603619966643
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.5831919999999999
$ws.Cells.Item(3, 14).Value = 1.749576
$ws.Cells.Item(3, 15).Value = 0.994689345773564
$ws.Cells.Item(3, 16).Value = 0.994689345773564
$ws.Cells.Item(3, 17).Value = 30.919087056444
$ws.Cells.Item(3, 18).Value = 185.514522338664
$ws.Cells.Item(3, 19).Value = 0.7402897018905847
$ws.Cells.Item(3, 20).Value = 0.700420809794138

# Row 4
$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 2).Value = 'Efnb2'
$ws.Cells.Item(4, 3).Value = 'Grm1'
$ws.Cells.Item(4, 4).Value = 'FAPs'
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.464838333333333
$ws.Cells.Item(4, 8).Value = 19.394515
$ws.Cells.Item(4, 9).Value = 0.09075212572810222
$ws.Cells.Item(4, 10).Value = 0.1287968964663749
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.003113666666666667
$ws.Cells.Item(4, 14).Value = 0.009341
$ws.Cells.Item(4, 15).Value = 0.005310654226435927
$ws.Cells.Item(4, 16).Value = 0.005310654226435927
$ws.Cells.Item(4, 17).Value = 0.02012935162388889
$ws.Cells.Item(4, 18).Value = 0.181164164615
$ws.Cells.Item(4, 19).Value = 0.0004819531600559907
$ws.Cells.Item(4, 20).Value = 0.0006839957825709844

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Efnb2'
$ws.Cells.Item(5, 3).Value = 'Grm1'
$ws.Cells.Item(5, 4).Value = 'Neutrophils'
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.464838333333333
$ws.Cells.Item(5, 8).Value = 19.394515
$ws.Cells.Item(5, 9).Value = 0.09075212572810222
$ws.Cells.Item(5, 10).Value = 0.1287968964663749
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.5831919999999999
$ws.Cells.Item(5, 14).Value = 1.749576
$ws.Cells.Item(5, 15).Value = 0.994689345773564
$ws.Cells.Item(5, 16).Value = 0.994689345773564
$ws.Cells.Item(5, 17).Value = 3.770241997293333
$ws.Cells.Item(5, 18).Value = 33.93217797563999
$ws.Cells.Item(5, 19).Value = 0.09027017256804623
$ws.Cells.Item(5, 20).Value = 0.1281129006838039

# Row 6
$ws.Cells.Item(6, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6, 2).Value = 'Efnb2'
$ws.Cells.Item(6, 3).Value = 'Grm1'
$ws.Cells.Item(6, 4).Value = 'FAPs'
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.6789063333333333
$ws.Cells.Item(6, 8).Value = 2.036719
$ws.Cells.Item(6, 9).Value = 0.00953035323444874
$ws.Cells.Item(6, 10).Value = 0.01352563269430035
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.003113666666666667
$ws.Cells.Item(6, 14).Value = 0.009341
$ws.Cells.Item(6, 15).Value = 0.005310654226435927
$ws.Cells.Item(6, 16).Value = 0.005310654226435927
$ws.Cells.Item(6, 17).Value = 0.002113888019888889
$ws.Cells.Item(6, 18).Value = 0.019024992179
$ws.Cells.Item(6, 19).Value = [double]"5.061241068395251E-05"
$ws.Cells.Item(6, 20).Value = [double]"7.182995843320612E-05"

# Row 7
$ws.Cells.Item(7, 1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(7, 2).Value = 'Efnb2'
$ws.Cells.Item(7, 3).Value = 'Grm1'
$ws.Cells.Item(7, 4).Value = 'Neutrophils'
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.6789063333333333
$ws.Cells.Item(7, 8).Value = 2.036719
$ws.Cells.Item(7, 9).Value = 0.00953035323444874
$ws.Cells.Item(7, 10).Value = 0.01352563269430035
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.5831919999999999
$ws.Cells.Item(7, 14).Value = 1.749576
$ws.Cells.Item(7, 15).Value = 0.994689345773564
$ws.Cells.Item(7, 16).Value = 0.994689345773564
$ws.Cells.Item(7, 17).Value = 0.3959327423493332
$ws.Cells.Item(7, 18).Value = 3.563394681143999
$ws.Cells.Item(7, 19).Value = 0.009479740823764788
$ws.Cells.Item(7, 20).Value = 0.01345380273586715

# Row 8
$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Efnb2'
$ws.Cells.Item(8, 3).Value = 'Grm1'
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 10.1095025
$ws.Cells.Item(8, 8).Value = 20.219005
$ws.Cells.Item(8, 9).Value = 0.1419152026119597
$ws.Cells.Item(8, 10).Value = 0.1342722462324073
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.003113666666666667
$ws.Cells.Item(8, 14).Value = 0.009341
$ws.Cells.Item(8, 15).Value = 0.005310654226435927
$ws.Cells.Item(8, 16).Value = 0.005310654226435927
$ws.Cells.Item(8, 17).Value = 0.03147762095083333
$ws.Cells.Item(8, 18).Value = 0.188865725705
$ws.Cells.Item(8, 19).Value = 0.0007536625705467145
$ws.Cells.Item(8, 20).Value = 0.0007130734719471793

# Row 9
$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Efnb2'
$ws.Cells.Item(9, 3).Value = 'Grm1'
$ws.Cells.Item(9, 4).Value = 'Neutrophils'
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 10.1095025
$ws.Cells.Item(9, 8).Value = 20.219005
$ws.Cells.Item(9, 9).Value = 0.1419152026119597
$ws.Cells.Item(9, 10).Value = 0.1342722462324073
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.5831919999999999
$ws.Cells.Item(9, 14).Value = 1.749576
$ws.Cells.Item(9, 15).Value = 0.994689345773564
$ws.Cells.Item(9, 16).Value = 0.994689345773564
$ws.Cells.Item(9, 17).Value = 5.895780981979999
$ws.Cells.Item(9, 18).Value = 35.37468589187999
$ws.Cells.Item(9, 19).Value = 0.141161540041413
$ws.Cells.Item(9, 20).Value = 0.1335591727604601

# Row 10
$ws.Cells.Item(10, 1).Value = 'Neutrophils'
$ws.Cells.Item(10, 2).Value = 'Efnb2'
$ws.Cells.Item(10, 3).Value = 'Grm1'
$ws.Cells.Item(10, 4).Value = 'FAPs'
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.5288903333333334
$ws.Cells.Item(10, 8).Value = 1.586671
$ws.Cells.Item(10, 9).Value = 0.007424458207958987
$ws.Cells.Item(10, 10).Value = 0.01053691213795238
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.003113666666666667
$ws.Cells.Item(10, 14).Value = 0.009341
$ws.Cells.Item(10, 15).Value = 0.005310654226435927
$ws.Cells.Item(10, 16).Value = 0.005310654226435927
$ws.Cells.Item(10, 17).Value = 0.001646788201222222
$ws.Cells.Item(10, 18).Value = 0.014821093811
$ws.Cells.Item(10, 19).Value = [double]"3.942873036109431E-05"
$ws.Cells.Item(10, 20).Value = [double]"5.595789697900084E-05"

# Row 11
$ws.Cells.Item(11, 1).Value = 'Neutrophils'
$ws.Cells.Item(11, 2).Value = 'Efnb2'
$ws.Cells.Item(11, 3).Value = 'Grm1'
$ws.Cells.Item(11, 4).Value = 'Neutrophils'
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.5288903333333334
$ws.Cells.Item(11, 8).Value = 1.586671
$ws.Cells.Item(11, 9).Value = 0.007424458207958987
$ws.Cells.Item(11, 10).Value = 0.01053691213795238
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.5831919999999999
$ws.Cells.Item(11, 14).Value = 1.749576
$ws.Cells.Item(11, 15).Value = 0.994689345773564
$ws.Cells.Item(11, 16).Value = 0.994689345773564
$ws.Cells.Item(11, 17).Value = 0.3084446112773333
$ws.Cells.Item(11, 18).Value = 2.776001501496
$ws.Cells.Item(11, 19).Value = 0.007385029477597892
$ws.Cells.Item(11, 20).Value = 0.01048095424097338

# Row 12
$ws.Cells.Item(12, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(12, 2).Value = 'Efnb2'
$ws.Cells.Item(12, 3).Value = 'Grm1'
$ws.Cells.Item(12, 4).Value = 'FAPs'
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.4370873333333334
$ws.Cells.Item(12, 8).Value = 1.311262
$ws.Cells.Item(12, 9).Value = 0.006135745796503949
$ws.Cells.Item(12, 10).Value = 0.008707950472300633
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.003113666666666667
$ws.Cells.Item(12, 14).Value = 0.009341
$ws.Cells.Item(12, 15).Value = 0.005310654226435927
$ws.Cells.Item(12, 16).Value = 0.005310654226435927
$ws.Cells.Item(12, 17).Value = 0.001360944260222222
$ws.Cells.Item(12, 18).Value = 0.012248498342
$ws.Cells.Item(12, 19).Value = [double]"3.258482434654017E-05"
$ws.Cells.Item(12, 20).Value = [double]"4.624491397931808E-05"

# Row 13
$ws.Cells.Item(13, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 2).Value = 'Efnb2'
$ws.Cells.Item(13, 3).Value = 'Grm1'
$ws.Cells.Item(13, 4).Value = 'Neutrophils'
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.4370873333333334
$ws.Cells.Item(13, 8).Value = 1.311262
$ws.Cells.Item(13, 9).Value = 0.006135745796503949
$ws.Cells.Item(13, 10).Value = 0.008707950472300633
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.5831919999999999
$ws.Cells.Item(13, 14).Value = 1.749576
$ws.Cells.Item(13, 15).Value = 0.994689345773564
$ws.Cells.Item(13, 16).Value = 0.994689345773564
$ws.Cells.Item(13, 17).Value = 0.2549058361013333
$ws.Cells.Item(13, 18).Value = 2.294152524912
$ws.Cells.Item(13, 19).Value = 0.006103160972157409
$ws.Cells.Item(13, 20).Value = 0.008661705558321315
